# Natmi following Dr Hou advice
# Updates ligand/receptor expressing-cell counts (E,K: 1 -> 3) and the
# downstream derived expression / specificity metrics recomputed from the
# new cell counts, for every data row (2-10) of the LR-pairs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ E=3; G=6.540008666666668;  H=19.620026; I=0.2365207520404831; J=0.2365207520404831; K=3; M=4.358929333333333; N=13.076788; O=0.315486934607403;  P=0.315486934607403;  Q=28.50743561738756; R=256.5669205564881;  S=0.07461920703228965; T=0.07461920703228966 }
    3  = @{ E=3; G=6.540008666666668;  H=19.620026; I=0.2365207520404831; J=0.2365207520404831; K=3; M=1.781382333333333; N=5.344147;  O=0.1289313977653647; P=0.1289313977653647; Q=11.65025589864689; R=104.852303087822;   S=0.03049495116109473; T=0.03049495116109473 }
    4  = @{ E=3; G=6.540008666666668;  H=19.620026; I=0.2365207520404831; J=0.2365207520404831; K=3; M=7.676201333333334; N=23.028604; O=0.5555816676272323; P=0.5555816676272323; Q=50.20242324707824; R=451.8218092237041;  S=0.1314065938470987;  T=0.1314065938470987 }
    5  = @{ E=3; G=14.67485766666667; H=44.024573; I=0.5307192311682535; J=0.5307192311682536; K=3; M=4.358929333333333; N=13.076788; O=0.315486934607403;  P=0.315486934607403;  Q=63.9666675457249;  R=575.700007911524;   S=0.16743498337847;    T=0.16743498337847 }
    6  = @{ E=3; G=14.67485766666667; H=44.024573; I=0.5307192311682535; J=0.5307192311682536; K=3; M=1.781382333333333; N=5.344147;  O=0.1289313977653647; P=0.1289313977653647; Q=26.14153219158122; R=235.273789724231;   S=0.06842637229548265; T=0.06842637229548267 }
    7  = @{ E=3; G=14.67485766666667; H=44.024573; I=0.5307192311682535; J=0.5307192311682536; K=3; M=7.676201333333334; N=23.028604; O=0.5555816676272323; P=0.5555816676272323; Q=112.6471619873436; R=1013.824457886092;  S=0.2948578754943009;  T=0.2948578754943009 }
    8  = @{ E=3; G=6.436020999999999; H=19.308063; I=0.2327600167912634; J=0.2327600167912634; K=3; M=4.358929333333333; N=13.076788; O=0.315486934607403;  P=0.315486934607403;  Q=28.05416072684933; R=252.487446541644;   S=0.07343274419664332; T=0.07343274419664333 }
    9  = @{ E=3; G=6.436020999999999; H=19.308063; I=0.2327600167912634; J=0.2327600167912634; K=3; M=1.781382333333333; N=5.344147;  O=0.1289313977653647; P=0.1289313977653647; Q=11.46501410636233; R=103.185126957261;   S=0.03001007430878735; T=0.03001007430878735 }
    10 = @{ E=3; G=6.436020999999999; H=19.308063; I=0.2327600167912634; J=0.2327600167912634; K=3; M=7.676201333333334; N=23.028604; O=0.5555816676272323; P=0.5555816676272323; Q=50.20242324707824; R=451.8218092237041;  S=0.1314065938470987;  T=0.1314065938470987 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
